$d = $word.ActiveDocument

# --- Update the date line ---
$d.Content.Find.Execute("2023-08-24 Thursday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2023-08-25 Friday", 2)

# --- Update the division problems in the table ---
# Each entry is Row, Column, New value (old values are addressed positionally
# since one value, "79÷5=", repeats twice in the original table).
$updates = @(
    @{Row=1;  Col=1; New="47÷8="},
    @{Row=1;  Col=2; New="36÷5="},
    @{Row=1;  Col=3; New="50÷2="},
    @{Row=1;  Col=4; New="83÷9="},
    @{Row=1;  Col=5; New="38÷9="},

    @{Row=5;  Col=1; New="66÷2="},
    @{Row=5;  Col=2; New="95÷7="},
    @{Row=5;  Col=3; New="49÷4="},
    @{Row=5;  Col=4; New="40÷3="},
    @{Row=5;  Col=5; New="94÷8="},

    @{Row=9;  Col=1; New="48÷4="},
    @{Row=9;  Col=2; New="14÷4="},
    @{Row=9;  Col=3; New="69÷7="},
    @{Row=9;  Col=4; New="81÷8="},
    @{Row=9;  Col=5; New="39÷2="},

    @{Row=13; Col=1; New="24÷9="},
    @{Row=13; Col=2; New="77÷6="},
    @{Row=13; Col=3; New="93÷3="},
    @{Row=13; Col=4; New="46÷5="},
    @{Row=13; Col=5; New="85÷7="},

    @{Row=17; Col=1; New="69÷2="},
    @{Row=17; Col=2; New="28÷9="},
    @{Row=17; Col=3; New="19÷5="},
    @{Row=17; Col=4; New="25÷3="},
    @{Row=17; Col=5; New="15÷4="}
)

foreach ($u in $updates) {
    $table = $d.Tables(1)
    $cell = $table.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}
